# Apply new TPM data to Wnt2-Fzd8.xlsx
# - Row 4 (Inflammatory-Mac) becomes what was row 5 (MuSCs) with recalculated values
# - Row 5 becomes what was row 6 (Resolving-Mac) with recalculated values
# - Row 6 is deleted; "Inflammatory-Mac" is removed from the label set entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs) - update M:T with new TPM-derived values
$ws.Range("M2").Value = 2.544438666666667
$ws.Range("N2").Value = 7.633316000000001
$ws.Range("O2").Value = 0.201325300207035
$ws.Range("P2").Value = 0.201325300207035
$ws.Range("Q2").Value = 0.1400416634822222
$ws.Range("R2").Value = 1.26037497134
$ws.Range("S2").Value = 0.201325300207035
$ws.Range("T2").Value = 0.201325300207035

# Row 3 (FAPs) - update O,P,S,T with new TPM-derived values
$ws.Range("O3").Value = 0.6969390273602759
$ws.Range("P3").Value = 0.696939027360276
$ws.Range("S3").Value = 0.6969390273602759
$ws.Range("T3").Value = 0.696939027360276

# Row 4 - now represents the former MuSCs row values
$ws.Range("D4").Value = "MuSCs"
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.273916333333333
$ws.Range("N4").Value = 3.821749
$ws.Range("O4").Value = 0.1007969229547075
$ws.Range("P4").Value = 0.1007969229547075
$ws.Range("Q4").Value = 0.07011423179277779
$ws.Range("R4").Value = 0.6310280861350001
$ws.Range("S4").Value = 0.1007969229547075
$ws.Range("T4").Value = 0.1007969229547075

# Row 5 - now represents the former Resolving-Mac row values
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01186433333333333
$ws.Range("N5").Value = 0.035593
$ws.Range("O5").Value = 0.0009387494779816524
$ws.Range("P5").Value = 0.0009387494779816526
$ws.Range("Q5").Value = 0.0006529931327777777
$ws.Range("R5").Value = 0.005876938195
$ws.Range("S5").Value = 0.0009387494779816524
$ws.Range("T5").Value = 0.0009387494779816526

# Delete row 6 entirely (old Resolving-Mac row is now folded into row 5)
$ws.Rows(6).Delete()
